$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION"
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# Delete column A entirely; this shifts the remaining columns (old B:F) left
# to become the new A:E
$ws.Columns("A").Delete()
